# "Controle de componentes do Processador" - update status column (Plan1)
# Several "A fazer" (to-do) items were completed and are now marked "Feito" (done).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")
$ws.Activate()

# Mark the finished components as "Feito" (I15 stays "A fazer" / remains pending)
$doneRows = @(5, 6, 7, 8, 9, 10, 11, 12, 14, 16)
foreach ($r in $doneRows) {
    $ws.Range("I$r").Value = "Feito"
}

# Move the active selection to I17 (just below the table), matching the saved
# cursor position recorded in the workbook.
$ws.Range("I17").Select()
